$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 900
$ws.Range("I33").Value = 1000
$ws.Range("J33").Value = 600
$ws.Range("K33").Value = 1000
$ws.Range("L33").Value = 600
$ws.Range("M33").Value = -771
$ws.Range("N33").Value = -1058
# Row 70
$ws.Range("H70").Value = 1050.125
$ws.Range("I70").Value = 1066.6666
$ws.Range("J70").Value = 1000.5
$ws.Range("K70").Value = 3199.9998
$ws.Range("L70").Value = 3001.5
$ws.Range("M70").Value = -2929.9998
$ws.Range("N70").Value = -3541.5
# Row 73
$ws.Range("H73").Value = 1050.125
$ws.Range("I73").Value = 1066.6666
$ws.Range("J73").Value = 1000.5
$ws.Range("K73").Value = 3199.9998
$ws.Range("L73").Value = 3001.5
$ws.Range("M73").Value = -2263.9998
$ws.Range("N73").Value = -4873.5
# Row 100
$ws.Range("H100").Value = 3251.25
$ws.Range("I100").Value = 2668.3333
$ws.Range("K100").Value = 2668.3333
$ws.Range("M100").Value = -2127.3333
# Row 129
$ws.Range("H129").Value = 357841.2
$ws.Range("J129").Value = 500816.5
$ws.Range("L129").Value = 1502449.5
$ws.Range("N129").Value = -1512449.5
# Row 141
$ws.Range("H141").Value = 3812.9167
$ws.Range("I141").Value = 2941.6667
$ws.Range("K141").Value = 8825.000100000001
$ws.Range("M141").Value = -3645.000100000001

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 97.5
$ws.Range("I4").Value = 95
$ws.Range("K4").Value = 95
$ws.Range("M4").Value = 21
# Row 32
$ws.Range("H32").Value = 2740.4878
$ws.Range("I32").Value = 2368.389
$ws.Range("J32").Value = 5419.6
$ws.Range("K32").Value = 2368.389
$ws.Range("L32").Value = 5419.6
$ws.Range("M32").Value = -2081.389
$ws.Range("N32").Value = -5993.6
# Row 45
$ws.Range("H45").Value = 3333.9429
$ws.Range("I45").Value = 2709.158
$ws.Range("K45").Value = 2709.158
$ws.Range("M45").Value = -2332.158
# Row 74
$ws.Range("H74").Value = 2862.795
$ws.Range("I74").Value = 3096
$ws.Range("J74").Value = 1796.7142
$ws.Range("K74").Value = 3096
$ws.Range("L74").Value = 1796.7142
$ws.Range("M74").Value = -2222
$ws.Range("N74").Value = -3544.7142
# Row 77
$ws.Range("H77").Value = 2862.795
$ws.Range("I77").Value = 3096
$ws.Range("J77").Value = 1796.7142
$ws.Range("K77").Value = 15480
$ws.Range("L77").Value = 8983.571
$ws.Range("M77").Value = -11112
$ws.Range("N77").Value = -17719.571
# Row 132
$ws.Range("H132").Value = 13772.878
$ws.Range("I132").Value = 1362.7142
$ws.Range("J132").Value = 86165.5
$ws.Range("K132").Value = 4088.1426
$ws.Range("L132").Value = 258496.5
$ws.Range("M132").Value = -1558.1426
$ws.Range("N132").Value = -263556.5

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 278.125
$ws.Range("I22").Value = 282.14285
$ws.Range("K22").Value = 282.14285
$ws.Range("M22").Value = -109.14285
# Row 105
$ws.Range("H105").Value = 2944473
$ws.Range("I105").Value = 3643
$ws.Range("J105").Value = 7145658.5
$ws.Range("K105").Value = 3643
$ws.Range("L105").Value = 7145658.5
$ws.Range("M105").Value = -1896
$ws.Range("N105").Value = -7149152.5
# Row 134
$ws.Range("H134").Value = 3051.1724
$ws.Range("I134").Value = 3088.7144
$ws.Range("K134").Value = 9266.143199999999
$ws.Range("M134").Value = -6731.143199999999

$ws = $wb.Worksheets.Item("CRP")
# Row 51
$ws.Range("H51").Value = 6100
$ws.Range("I51").Value = 6100
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 6100
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -5364
$ws.Range("N51").ClearContents()
# Row 61
$ws.Range("H61").Value = 6100
$ws.Range("I61").Value = 6100
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 6100
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -5752
$ws.Range("N61").ClearContents()
# Row 62
$ws.Range("H62").Value = 4557.8887
$ws.Range("I62").Value = 2751.25
$ws.Range("J62").Value = 6003.2
$ws.Range("K62").Value = 2751.25
$ws.Range("L62").Value = 6003.2
$ws.Range("M62").Value = -2127.25
$ws.Range("N62").Value = -7251.2
# Row 65
$ws.Range("H65").Value = 4557.8887
$ws.Range("I65").Value = 2751.25
$ws.Range("J65").Value = 6003.2
$ws.Range("K65").Value = 13756.25
$ws.Range("L65").Value = 30016
$ws.Range("M65").Value = -10636.25
$ws.Range("N65").Value = -36256
# Row 68
$ws.Range("H68").Value = 67467.5
$ws.Range("J68").Value = 67467.5
$ws.Range("L68").Value = 67467.5
$ws.Range("N68").Value = -68965.5
# Row 71
$ws.Range("H71").Value = 67467.5
$ws.Range("J71").Value = 67467.5
$ws.Range("L71").Value = 202402.5
$ws.Range("N71").Value = -209890.5
# Row 99
$ws.Range("H99").Value = 20837398
$ws.Range("I99").Value = 4905197
$ws.Range("J99").Value = 45459892
$ws.Range("K99").Value = 4905197
$ws.Range("L99").Value = 45459892
$ws.Range("M99").Value = -4903699
$ws.Range("N99").Value = -45462888
# Row 105
$ws.Range("I105").Value = 20833882
$ws.Range("J105").Value = 389
$ws.Range("K105").Value = 20833882
$ws.Range("L105").Value = 389
$ws.Range("M105").Value = -20832135
$ws.Range("N105").Value = -3883
# Row 126
$ws.Range("H126").Value = 20837398
$ws.Range("I126").Value = 4905197
$ws.Range("J126").Value = 45459892
$ws.Range("K126").Value = 14715591
$ws.Range("L126").Value = 136379676
$ws.Range("M126").Value = -14713121
$ws.Range("N126").Value = -136384616

$ws = $wb.Worksheets.Item("CUL")
# Row 118
$ws.Range("H118").Value = 125002950
$ws.Range("I118").Value = 500000000
$ws.Range("J118").Value = 3933.3333
$ws.Range("K118").Value = 1500000000
$ws.Range("L118").Value = 11799.9999
$ws.Range("M118").Value = -1499998757
$ws.Range("N118").Value = -14285.9999
# Row 131
$ws.Range("H131").Value = 101817.445
$ws.Range("J131").Value = 103896.16
$ws.Range("L131").Value = 311688.48
$ws.Range("N131").Value = -321768.48

$ws = $wb.Worksheets.Item("GSM")
# Row 58
$ws.Range("H58").Value = 25005624
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 25005624
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 25005624
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -25006178
# Row 80
$ws.Range("H80").Value = 3916.5334
$ws.Range("I80").Value = 3296
$ws.Range("J80").Value = 4625.7144
$ws.Range("K80").Value = 3296
$ws.Range("L80").Value = 4625.7144
$ws.Range("M80").Value = -2298
$ws.Range("N80").Value = -6621.7144
# Row 83
$ws.Range("H83").Value = 3916.5334
$ws.Range("I83").Value = 3296
$ws.Range("J83").Value = 4625.7144
$ws.Range("K83").Value = 16480
$ws.Range("L83").Value = 23128.572
$ws.Range("M83").Value = -11488
$ws.Range("N83").Value = -33112.572

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 2819.1667
$ws.Range("I61").Value = 1242.5
$ws.Range("J61").Value = 7549.1665
$ws.Range("K61").Value = 1242.5
$ws.Range("L61").Value = 7549.1665
$ws.Range("M61").Value = -1040.5
$ws.Range("N61").Value = -7953.1665
# Row 94
$ws.Range("H94").Value = 15000
$ws.Range("J94").Value = 15000
$ws.Range("L94").Value = 15000
$ws.Range("N94").Value = -16352
# Row 113
$ws.Range("H113").Value = 2819.1667
$ws.Range("I113").Value = 1242.5
$ws.Range("J113").Value = 7549.1665
$ws.Range("K113").Value = 1242.5
$ws.Range("L113").Value = 7549.1665
$ws.Range("M113").Value = 927.5
$ws.Range("N113").Value = -11889.1665

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 2675083.8
$ws.Range("J107").Value = 4133684.2
$ws.Range("L107").Value = 12401052.6
$ws.Range("N107").Value = -12404892.6

Write-Host "Applied Typhon_Profits updates"